# "Nueva tarea en el productBacklog"
# Insert a new Product Backlog Item row right after the first existing
# item ("Preparacion del Entorno de desarrollo"), pushing the rest of
# the backlog down by one row, and renumber the Id column to stay
# sequential (1..7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Insert a new row at row 4 (shifts old rows 4-9 down to 5-10, etc.)
$ws.Rows.Item(4).Insert()

# Keep the same row height / look as the rest of the table
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(3).RowHeight

# New backlog item content
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Ejecutar un ejemplo de Prueba e instalar en apk"

# Renumber the Id column for the items that were shifted down
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

# Match formatting of the surrounding rows for the new/shifted cells
$ws.Range("C4:E4").Style = $ws.Range("C5").Style
$ws.Range("B4").Style = $ws.Range("B5").Style

$ws.Range("A4").Select()
